# Add 4 new data rows (41-44) to the GCF_File_Usage sheet, matching the
# same styling/number format as the existing data rows (e.g. row 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: Col A is a date/time serial value, columns B:O are plain numbers.
$newRows = @(
    @{ Row = 41; A = 45779.894999999997; Values = @(10,6,360,537,512,565,4214,565,2728,268,510,30,4843,6400) },
    @{ Row = 42; A = 45782.992569444446; Values = @(10,6,360,538,513,572,4292,572,2842,278,519,30,4875,6448) },
    @{ Row = 43; A = 45783.482037037036; Values = @(10,6,360,538,513,577,4305,577,3084,284,527,30,4875,6448) },
    @{ Row = 44; A = 45783.483414351853; Values = @(10,6,360,538,513,577,4305,577,3084,284,527,30,4920,6448) }
)

# Template row (last existing data row) to copy styles from.
$templateRow = 40

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Copy formatting from the template row for columns A:O so the new rows
    # look the same as the existing data (same number format / style ids).
    $srcRange = $ws.Range("A$templateRow`:O$templateRow")
    $dstRange = $ws.Range("A$r`:O$r")
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    # Column A: date/time serial value
    $ws.Cells.Item($r, 1).Value = $entry.A

    # Columns B..O (2..15): plain numeric values
    $col = 2
    foreach ($v in $entry.Values) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}

$excel.CutCopyMode = $false
